# Auto-generated Excel COM-interop script
# Applies numeric cell-value corrections to the Leve profit tables
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- Cell value updates (Sheet, CellRef, NewValue) ---
$updates = @(
    @{Sheet="ALC"; Cell="H33"; Value=206.88889},
    @{Sheet="ALC"; Cell="I33"; Value=148.85715},
    @{Sheet="ALC"; Cell="J33"; Value=410},
    @{Sheet="ALC"; Cell="K33"; Value=148.85715},
    @{Sheet="ALC"; Cell="L33"; Value=410},
    @{Sheet="ALC"; Cell="M33"; Value=80.14285000000001},
    @{Sheet="ALC"; Cell="N33"; Value=-868},
    @{Sheet="ALC"; Cell="H112"; Value=2078.15},
    @{Sheet="ALC"; Cell="J112"; Value=2341.4375},
    @{Sheet="ALC"; Cell="L112"; Value=7024.3125},
    @{Sheet="ALC"; Cell="N112"; Value=-9240.3125},
    @{Sheet="ALC"; Cell="H129"; Value=971.1667},
    @{Sheet="ALC"; Cell="J129"; Value=1004.8222},
    @{Sheet="ALC"; Cell="L129"; Value=3014.4666},
    @{Sheet="ALC"; Cell="N129"; Value=-13014.4666},
    @{Sheet="ALC"; Cell="H137"; Value=5429.788},
    @{Sheet="ALC"; Cell="I137"; Value=5726.2104},
    @{Sheet="ALC"; Cell="J137"; Value=5027.5},
    @{Sheet="ALC"; Cell="K137"; Value=17178.6312},
    @{Sheet="ALC"; Cell="L137"; Value=15082.5},
    @{Sheet="ALC"; Cell="M137"; Value=-14628.6312},
    @{Sheet="ALC"; Cell="N137"; Value=-20182.5},
    @{Sheet="ARM"; Cell="H32"; Value=15389.598},
    @{Sheet="ARM"; Cell="I32"; Value=10456.973},
    @{Sheet="ARM"; Cell="J32"; Value=20322.223},
    @{Sheet="ARM"; Cell="K32"; Value=10456.973},
    @{Sheet="ARM"; Cell="L32"; Value=20322.223},
    @{Sheet="ARM"; Cell="M32"; Value=-10169.973},
    @{Sheet="ARM"; Cell="N32"; Value=-20896.223},
    @{Sheet="ARM"; Cell="H34"; Value=26000},
    @{Sheet="ARM"; Cell="J34"; Value=26000},
    @{Sheet="ARM"; Cell="L34"; Value=26000},
    @{Sheet="ARM"; Cell="N34"; Value=-26542},
    @{Sheet="ARM"; Cell="H45"; Value=1068.579},
    @{Sheet="ARM"; Cell="I45"; Value=1099.1538},
    @{Sheet="ARM"; Cell="J45"; Value=1002.3333},
    @{Sheet="ARM"; Cell="K45"; Value=1099.1538},
    @{Sheet="ARM"; Cell="L45"; Value=1002.3333},
    @{Sheet="ARM"; Cell="M45"; Value=-722.1538},
    @{Sheet="ARM"; Cell="N45"; Value=-1756.3333},
    @{Sheet="ARM"; Cell="H61"; Value=2375.6667},
    @{Sheet="ARM"; Cell="I61"; Value=1914.4286},
    @{Sheet="ARM"; Cell="K61"; Value=1914.4286},
    @{Sheet="ARM"; Cell="M61"; Value=-1702.4286},
    @{Sheet="ARM"; Cell="H136"; Value=2375.6667},
    @{Sheet="ARM"; Cell="I136"; Value=1914.4286},
    @{Sheet="ARM"; Cell="K136"; Value=5743.2858},
    @{Sheet="ARM"; Cell="M136"; Value=-3193.2858},
    @{Sheet="BSM"; Cell="H97"; Value=22106.467},
    @{Sheet="BSM"; Cell="I97"; Value=4026.2856},
    @{Sheet="BSM"; Cell="J97"; Value=37926.625},
    @{Sheet="BSM"; Cell="K97"; Value=4026.2856},
    @{Sheet="BSM"; Cell="L97"; Value=37926.625},
    @{Sheet="BSM"; Cell="M97"; Value=-3035.2856},
    @{Sheet="BSM"; Cell="N97"; Value=-39908.625},
    @{Sheet="BSM"; Cell="H118"; Value=28888.572},
    @{Sheet="BSM"; Cell="J118"; Value=28888.572},
    @{Sheet="BSM"; Cell="L118"; Value=28888.572},
    @{Sheet="BSM"; Cell="N118"; Value=-32202.572},
    @{Sheet="CRP"; Cell="H31"; Value=5359.2188},
    @{Sheet="CRP"; Cell="I31"; Value=2950},
    @{Sheet="CRP"; Cell="J31"; Value=5703.393},
    @{Sheet="CRP"; Cell="K31"; Value=2950},
    @{Sheet="CRP"; Cell="L31"; Value=5703.393},
    @{Sheet="CRP"; Cell="M31"; Value=-2655},
    @{Sheet="CRP"; Cell="N31"; Value=-6293.393},
    @{Sheet="CRP"; Cell="H34"; Value=5359.2188},
    @{Sheet="CRP"; Cell="I34"; Value=2950},
    @{Sheet="CRP"; Cell="J34"; Value=5703.393},
    @{Sheet="CRP"; Cell="K34"; Value=2950},
    @{Sheet="CRP"; Cell="L34"; Value=5703.393},
    @{Sheet="CRP"; Cell="M34"; Value=-2748},
    @{Sheet="CRP"; Cell="N34"; Value=-6107.393},
    @{Sheet="CRP"; Cell="H58"; Value=2277.6936},
    @{Sheet="CRP"; Cell="I58"; Value=1813.0546},
    @{Sheet="CRP"; Cell="J58"; Value=5928.4287},
    @{Sheet="CRP"; Cell="K58"; Value=1813.0546},
    @{Sheet="CRP"; Cell="L58"; Value=5928.4287},
    @{Sheet="CRP"; Cell="M58"; Value=-1610.0546},
    @{Sheet="CRP"; Cell="N58"; Value=-6334.4287},
    @{Sheet="CRP"; Cell="H62"; Value=62505708},
    @{Sheet="CRP"; Cell="J62"; Value=5945},
    @{Sheet="CRP"; Cell="L62"; Value=5945},
    @{Sheet="CRP"; Cell="N62"; Value=-7193},
    @{Sheet="CRP"; Cell="H65"; Value=62505708},
    @{Sheet="CRP"; Cell="J65"; Value=5945},
    @{Sheet="CRP"; Cell="L65"; Value=29725},
    @{Sheet="CRP"; Cell="N65"; Value=-35965},
    @{Sheet="CRP"; Cell="H132"; Value=2723},
    @{Sheet="CRP"; Cell="I132"; Value=2123.7317},
    @{Sheet="CRP"; Cell="J132"; Value=5453},
    @{Sheet="CRP"; Cell="K132"; Value=6371.195099999999},
    @{Sheet="CRP"; Cell="L132"; Value=16359},
    @{Sheet="CRP"; Cell="M132"; Value=-3841.195099999999},
    @{Sheet="CRP"; Cell="N132"; Value=-21419},
    @{Sheet="CRP"; Cell="H136"; Value=2277.6936},
    @{Sheet="CRP"; Cell="I136"; Value=1813.0546},
    @{Sheet="CRP"; Cell="J136"; Value=5928.4287},
    @{Sheet="CRP"; Cell="K136"; Value=5439.1638},
    @{Sheet="CRP"; Cell="L136"; Value=17785.2861},
    @{Sheet="CRP"; Cell="M136"; Value=-2889.1638},
    @{Sheet="CRP"; Cell="N136"; Value=-22885.2861},
    @{Sheet="CRP"; Cell="H138"; Value=45118.188},
    @{Sheet="CRP"; Cell="J138"; Value=45118.188},
    @{Sheet="CRP"; Cell="L138"; Value=45118.188},
    @{Sheet="CRP"; Cell="N138"; Value=-55398.188},
    @{Sheet="CRP"; Cell="H139"; Value=112543.336},
    @{Sheet="CRP"; Cell="J139"; Value=112543.336},
    @{Sheet="CRP"; Cell="L139"; Value=112543.336},
    @{Sheet="CRP"; Cell="N139"; Value=-122823.336},
    @{Sheet="CRP"; Cell="H140"; Value=115520},
    @{Sheet="CRP"; Cell="J140"; Value=115520},
    @{Sheet="CRP"; Cell="L140"; Value=115520},
    @{Sheet="CRP"; Cell="N140"; Value=-125880},
    @{Sheet="CRP"; Cell="H141"; Value=28396.969},
    @{Sheet="CRP"; Cell="J141"; Value=28396.969},
    @{Sheet="CRP"; Cell="L141"; Value=28396.969},
    @{Sheet="CRP"; Cell="N141"; Value=-38756.969},
    @{Sheet="CUL"; Cell="H5"; Value=1971.3334},
    @{Sheet="CUL"; Cell="J5"; Value=3878.2},
    @{Sheet="CUL"; Cell="L5"; Value=11634.6},
    @{Sheet="CUL"; Cell="N5"; Value=-11858.6},
    @{Sheet="CUL"; Cell="H105"; Value=4993.3335},
    @{Sheet="CUL"; Cell="J105"; Value=4993.3335},
    @{Sheet="CUL"; Cell="L105"; Value=14980.0005},
    @{Sheet="CUL"; Cell="N105"; Value=-20222.0005},
    @{Sheet="CUL"; Cell="H122"; Value=2565.7031},
    @{Sheet="CUL"; Cell="I122"; Value=789.7},
    @{Sheet="CUL"; Cell="J122"; Value=2894.5925},
    @{Sheet="CUL"; Cell="K122"; Value=7107.3},
    @{Sheet="CUL"; Cell="L122"; Value=26051.3325},
    @{Sheet="CUL"; Cell="M122"; Value=-4657.3},
    @{Sheet="CUL"; Cell="N122"; Value=-30951.3325},
    @{Sheet="CUL"; Cell="H131"; Value=7816610},
    @{Sheet="CUL"; Cell="I131"; Value=33347652},
    @{Sheet="CUL"; Cell="J131"; Value=984.7755},
    @{Sheet="CUL"; Cell="K131"; Value=100042956},
    @{Sheet="CUL"; Cell="L131"; Value=2954.3265},
    @{Sheet="CUL"; Cell="M131"; Value=-100037916},
    @{Sheet="CUL"; Cell="N131"; Value=-13034.3265},
    @{Sheet="CUL"; Cell="H132"; Value=1830.4722},
    @{Sheet="CUL"; Cell="J132"; Value=2213.261},
    @{Sheet="CUL"; Cell="L132"; Value=19919.349},
    @{Sheet="CUL"; Cell="N132"; Value=-24979.349},
    @{Sheet="CUL"; Cell="H133"; Value=3462.1738},
    @{Sheet="CUL"; Cell="I133"; Value=3359.2307},
    @{Sheet="CUL"; Cell="J133"; Value=3596},
    @{Sheet="CUL"; Cell="K133"; Value=10077.6921},
    @{Sheet="CUL"; Cell="L133"; Value=10788},
    @{Sheet="CUL"; Cell="M133"; Value=-5017.6921},
    @{Sheet="CUL"; Cell="N133"; Value=-20908},
    @{Sheet="CUL"; Cell="H135"; Value=1971.3334},
    @{Sheet="CUL"; Cell="J135"; Value=3878.2},
    @{Sheet="CUL"; Cell="L135"; Value=34903.8},
    @{Sheet="CUL"; Cell="N135"; Value=-39973.8},
    @{Sheet="CUL"; Cell="H137"; Value=9148.9},
    @{Sheet="CUL"; Cell="I137"; Value=3158.6667},
    @{Sheet="CUL"; Cell="K137"; Value=9476.000100000001},
    @{Sheet="CUL"; Cell="M137"; Value=-4376.000100000001},
    @{Sheet="GSM"; Cell="H70"; Value=5965},
    @{Sheet="GSM"; Cell="I70"; Value=5395.346},
    @{Sheet="GSM"; Cell="K70"; Value=5395.346},
    @{Sheet="GSM"; Cell="M70"; Value=-5125.346},
    @{Sheet="GSM"; Cell="H73"; Value=5965},
    @{Sheet="GSM"; Cell="I73"; Value=5395.346},
    @{Sheet="GSM"; Cell="K73"; Value=5395.346},
    @{Sheet="GSM"; Cell="M73"; Value=-4459.346},
    @{Sheet="GSM"; Cell="H117"; Value=27903.334},
    @{Sheet="GSM"; Cell="J117"; Value=27903.334},
    @{Sheet="GSM"; Cell="L117"; Value=27903.334},
    @{Sheet="GSM"; Cell="N117"; Value=-34787.334},
    @{Sheet="GSM"; Cell="H119"; Value=39766.668},
    @{Sheet="GSM"; Cell="J119"; Value=39766.668},
    @{Sheet="GSM"; Cell="L119"; Value=39766.668},
    @{Sheet="GSM"; Cell="N119"; Value=-49442.668},
    @{Sheet="GSM"; Cell="H132"; Value=5239.75},
    @{Sheet="GSM"; Cell="I132"; Value=0},
    @{Sheet="GSM"; Cell="J132"; Value=5239.75},
    @{Sheet="GSM"; Cell="K132"; Value=0},
    @{Sheet="GSM"; Cell="L132"; Value=15719.25},
    @{Sheet="GSM"; Cell="N132"; Value=-20779.25},
    @{Sheet="LTW"; Cell="H56"; Value=12778.25},
    @{Sheet="LTW"; Cell="I56"; Value=5500},
    @{Sheet="LTW"; Cell="J56"; Value=20056.5},
    @{Sheet="LTW"; Cell="K56"; Value=5500},
    @{Sheet="LTW"; Cell="L56"; Value=20056.5},
    @{Sheet="LTW"; Cell="M56"; Value=-4809},
    @{Sheet="LTW"; Cell="N56"; Value=-21438.5},
    @{Sheet="LTW"; Cell="H68"; Value=1971.875},
    @{Sheet="LTW"; Cell="I68"; Value=595},
    @{Sheet="LTW"; Cell="J68"; Value=2168.5715},
    @{Sheet="LTW"; Cell="K68"; Value=595},
    @{Sheet="LTW"; Cell="L68"; Value=2168.5715},
    @{Sheet="LTW"; Cell="M68"; Value=154},
    @{Sheet="LTW"; Cell="N68"; Value=-3666.5715},
    @{Sheet="LTW"; Cell="H71"; Value=1971.875},
    @{Sheet="LTW"; Cell="I71"; Value=595},
    @{Sheet="LTW"; Cell="J71"; Value=2168.5715},
    @{Sheet="LTW"; Cell="K71"; Value=2975},
    @{Sheet="LTW"; Cell="L71"; Value=10842.8575},
    @{Sheet="LTW"; Cell="M71"; Value=769},
    @{Sheet="LTW"; Cell="N71"; Value=-18330.8575},
    @{Sheet="LTW"; Cell="H115"; Value=30000},
    @{Sheet="LTW"; Cell="I115"; Value=0},
    @{Sheet="LTW"; Cell="J115"; Value=30000},
    @{Sheet="LTW"; Cell="K115"; Value=0},
    @{Sheet="LTW"; Cell="L115"; Value=30000},
    @{Sheet="LTW"; Cell="N115"; Value=-32350},
    @{Sheet="LTW"; Cell="H128"; Value=42780},
    @{Sheet="LTW"; Cell="J128"; Value=42780},
    @{Sheet="LTW"; Cell="L128"; Value=42780},
    @{Sheet="LTW"; Cell="N128"; Value=-52740},
    @{Sheet="WVR"; Cell="H136"; Value=4858.316},
    @{Sheet="WVR"; Cell="I136"; Value=1564.2727},
    @{Sheet="WVR"; Cell="J136"; Value=9387.625},
    @{Sheet="WVR"; Cell="K136"; Value=4692.8181},
    @{Sheet="WVR"; Cell="L136"; Value=28162.875},
    @{Sheet="WVR"; Cell="M136"; Value=-2142.8181},
    @{Sheet="WVR"; Cell="N136"; Value=-33262.875}
)

# --- Cells whose values are removed entirely (no longer present) ---
$deletions = @(
    @{Sheet="GSM"; Cell="M132"},
    @{Sheet="LTW"; Cell="M115"}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

foreach ($d in $deletions) {
    $ws = $wb.Worksheets.Item($d.Sheet)
    $ws.Range($d.Cell).ClearContents()
}

Write-Host "Applied $($updates.Count) cell updates and $($deletions.Count) deletions."
